$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 5, shifting existing rows 5..78 down to 6..79
$ws.Rows("5:5").Insert()

# Populate the new row 5 with the new data record
$ws.Range("A5").Value = 11
$ws.Range("B5").Value = "Vega Monumental Concepción"
$ws.Range("C5").Value = "Bíobío"
$ws.Range("D5").Value = 44860
$ws.Range("E5").Value = 8
$ws.Range("F5").Value = 100112013
$ws.Range("G5").Value = "Alcachofa"
$ws.Range("H5").Value = "Española"
$ws.Range("I5").Value = "Segunda"
$ws.Range("J5").Value = 180
$ws.Range("K5").Value = 8000
$ws.Range("L5").Value = 8500
$ws.Range("M5").Value = 8222
$ws.Range("N5").Value = "$/caja 40 unidades"
$ws.Range("O5").Value = "Provincia de Limarí"
$ws.Range("P5").Value = 206
$ws.Range("Q5").Value = 40
$ws.Range("R5").Value = "Hortaliza"
